# Round the numeric data in B2:E13 to the nearest integer, in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:E13")

for ($r = 1; $r -le $range.Rows.Count; $r++) {
    for ($c = 1; $c -le $range.Columns.Count; $c++) {
        $cell = $range.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $cell.Value2 = [Math]::Floor([double]$val + 0.5)
        }
    }
}
